$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.146258333333333
$ws.Range("H2").Value = 6.438775
$ws.Range("I2").Value = 0.3745961416936294
$ws.Range("J2").Value = 0.3745961416936293
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.553279333333334
$ws.Range("N2").Value = 7.659838000000001
$ws.Range("O2").Value = 0.1645043904057808
$ws.Range("P2").Value = 0.1645043904057808
$ws.Range("Q2").Value = 5.479997046494446
$ws.Range("R2").Value = 49.31997341845
$ws.Range("S2").Value = 0.06162270993766799
$ws.Range("T2").Value = 0.06162270993766797

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.146258333333333
$ws.Range("H3").Value = 6.438775
$ws.Range("I3").Value = 0.3745961416936294
$ws.Range("J3").Value = 0.3745961416936293
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.058662
$ws.Range("N3").Value = 24.175986
$ws.Range("O3").Value = 0.5192088709172035
$ws.Range("P3").Value = 0.5192088709172035
$ws.Range("Q3").Value = 17.29597047301667
$ws.Range("R3").Value = 155.66373425715
$ws.Range("S3").Value = 0.1944936397786901
$ws.Range("T3").Value = 0.1944936397786901

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.146258333333333
$ws.Range("H4").Value = 6.438775
$ws.Range("I4").Value = 0.3745961416936294
$ws.Range("J4").Value = 0.3745961416936293
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.909099333333334
$ws.Range("N4").Value = 14.727298
$ws.Range("O4").Value = 0.3162867386770157
$ws.Range("P4").Value = 0.3162867386770157
$ws.Range("Q4").Value = 10.53619535332778
$ws.Range("R4").Value = 94.82575817995
$ws.Range("S4").Value = 0.1184797919772713
$ws.Range("T4").Value = 0.1184797919772713

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.966337333333333
$ws.Range("H5").Value = 5.899012
$ws.Range("I5").Value = 0.3431937185263377
$ws.Range("J5").Value = 0.3431937185263377
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.553279333333334
$ws.Range("N5").Value = 7.659838000000001
$ws.Range("O5").Value = 0.1645043904057808
$ws.Range("P5").Value = 0.1645043904057808
$ws.Range("Q5").Value = 5.020608475561779
$ws.Range("R5").Value = 45.185476280056
$ws.Range("S5").Value = 0.0564568734572683
$ws.Range("T5").Value = 0.05645687345726829

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.966337333333333
$ws.Range("H6").Value = 5.899012
$ws.Range("I6").Value = 0.3431937185263377
$ws.Range("J6").Value = 0.3431937185263377
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.058662
$ws.Range("N6").Value = 24.175986
$ws.Range("O6").Value = 0.5192088709172035
$ws.Range("P6").Value = 0.5192088709172035
$ws.Range("Q6").Value = 15.84604794731467
$ws.Range("R6").Value = 142.614431525832
$ws.Range("S6").Value = 0.1781892231019363
$ws.Range("T6").Value = 0.1781892231019363

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.966337333333333
$ws.Range("H7").Value = 5.899012
$ws.Range("I7").Value = 0.3431937185263377
$ws.Range("J7").Value = 0.3431937185263377
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.909099333333334
$ws.Range("N7").Value = 14.727298
$ws.Range("O7").Value = 0.3162867386770157
$ws.Range("P7").Value = 0.3162867386770157
$ws.Range("Q7").Value = 9.652945292175112
$ws.Range("R7").Value = 86.87650762957601
$ws.Range("S7").Value = 0.108547621967133
$ws.Range("T7").Value = 0.108547621967133

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.616930333333333
$ws.Range("H8").Value = 4.850791
$ws.Range("I8").Value = 0.282210139780033
$ws.Range("J8").Value = 0.282210139780033
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.553279333333334
$ws.Range("N8").Value = 7.659838000000001
$ws.Range("O8").Value = 0.1645043904057808
$ws.Range("P8").Value = 0.1645043904057808
$ws.Range("Q8").Value = 4.128474803539778
$ws.Range("R8").Value = 37.156273231858
$ws.Range("S8").Value = 0.04642480701084452
$ws.Range("T8").Value = 0.04642480701084451

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.616930333333333
$ws.Range("H9").Value = 4.850791
$ws.Range("I9").Value = 0.282210139780033
$ws.Range("J9").Value = 0.282210139780033
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.058662
$ws.Range("N9").Value = 24.175986
$ws.Range("O9").Value = 0.5192088709172035
$ws.Range("P9").Value = 0.5192088709172035
$ws.Range("Q9").Value = 13.03029503388067
$ws.Range("R9").Value = 117.272655304926
$ws.Range("S9").Value = 0.1465260080365771
$ws.Range("T9").Value = 0.1465260080365771

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.616930333333333
$ws.Range("H10").Value = 4.850791
$ws.Range("I10").Value = 0.282210139780033
$ws.Range("J10").Value = 0.282210139780033
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.909099333333334
$ws.Range("N10").Value = 14.727298
$ws.Range("O10").Value = 0.3162867386770157
$ws.Range("P10").Value = 0.3162867386770157
$ws.Range("Q10").Value = 7.937671621413112
$ws.Range("R10").Value = 71.439044592718
$ws.Range("S10").Value = 0.08925932473261136
$ws.Range("T10").Value = 0.08925932473261136

